$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Add new row 14 data
$ws.Range("A14").Value = "Kun jij dit afhandelen?"
$ws.Range("B14").Value = "mailmind.test@zohomail.eu"
$ws.Range("C14").Value = "Testmail #3: Kun jij dit afhandelen?"
$ws.Range("D14").Value = "Planning / Afspraak"
$ws.Range("E14").Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$ws.Range("F14").Value = "2025-08-03 14:38:07"
$ws.Range("G14").Value = "Ja"
$ws.Range("H14").Value = "Ja"
$ws.Range("I14").Value = "Nee"
$ws.Range("J14").Value = "Nee"

# Update Dashboard sheet count
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B3").Value = 4

# Extend conditional formatting ranges to include row 14
foreach ($col in @("D", "G", "H", "I", "J")) {
    $rng = $ws.Range("${col}2")
    $fcs = $rng.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($ws.Range("${col}2:${col}14"))
    }
}
